$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.664.27'
$ws.Range('E2').Value = '  -3.92%  '
$ws.Range('D3').Value = '3.690.76'
$ws.Range('E3').Value = '  -4.76%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.70%  '
$ws.Range('D7').Value = '3.687.31'
$ws.Range('E7').Value = '  -4.70%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.631'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.716'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.161'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000294'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.48%  '
$ws.Range('D15').Value = '4.267.15'
$ws.Range('E15').Value = '  -5.05%  '
$ws.Range('D16').Value = '3.680.95'
$ws.Range('E16').Value = '  -5.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.45%  '
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.28%  '
$ws.Range('D21').Value = '68.278.30'
$ws.Range('E21').Value = '  -4.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '412.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.34%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.55%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.63'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.97'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -15.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.47'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.24%  '
$ws.Range('E34').Value = '  -5.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.89%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '43.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.33%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '601.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.93%  '
$ws.Range('D38').Value = '0.0₃0882'
$ws.Range('E38').Value = '  -11.97%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.403'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.48%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -6.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0442'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.65'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.74'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.60%  '
$ws.Range('E48').Value = '  -6.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.82%  '
$ws.Range('D50').Value = '2.724.15'
$ws.Range('E50').Value = '  -6.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.66%  '
